$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# 2. Replace the Gen generation-count values in column A with MaxFES fractions
$newA = @{
    2  = 0
    3  = 0.001
    4  = 0.01
    5  = 0.1
    6  = 0.2
    7  = 0.3
    8  = 0.4
    9  = 0.5
    10 = 0.6
    11 = 0.7
    12 = 0.8
    13 = 0.9
    14 = 1
}
foreach ($r in $newA.Keys) {
    $ws.Cells.Item($r, 1).Value = $newA[$r]
}

# 3. Replace column AZ ("Run 50") data with the new recomputed Mean values
$newMean = @{
    2  = 527994398.3901558
    3  = 318766554.5585039
    4  = 31547127.0260138
    5  = 1945917.49501263
    6  = 707688.3074404
    7  = 357758.85186398
    8  = 219975.62463317
    9  = 148395.16356122
    10 = 106460.94171288
    11 = 81583.26439881
    12 = 56636.39948255
    13 = 42869.45758638
    14 = 33838.83570655
}
foreach ($r in $newMean.Keys) {
    $ws.Cells.Item($r, 52).Value = $newMean[$r]
}

# 4. AZ1 now holds the "Mean" header (previously held by BA1)
$ws.Cells.Item(1, 52).Value = "Mean"

# 5. Delete the now-redundant last column (BA), which held the old "Run 50"
#    header slot data before step 3/4 overwrote AZ; after the overwrite BA
#    is duplicate/obsolete and must be removed so the sheet shrinks to A:AZ.
$ws.Columns.Item(53).Delete()

Write-Host "Done editing sheet"
